$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, pushing existing rows 124:192 down to 125:193
$ws.Rows(124).Insert()

# Populate the new row 124 with the new data record
$ws.Range("A124").Value = 7
$ws.Range("B124").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C124").Value = "Ñuble"
$ws.Range("D124").Value = 44875
$ws.Range("E124").Value = 16
$ws.Range("F124").Value = 100112028
$ws.Range("G124").Value = "Sandia"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 300
$ws.Range("K124").Value = 800
$ws.Range("L124").Value = 900
$ws.Range("M124").Value = 850
$ws.Range("N124").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O124").Value = "Perú"
$ws.Range("P124").Value = 850
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = "Hortaliza"

# Match the date format used by the other Fecha cells in column D
$ws.Range("D124").NumberFormat = $ws.Range("D125").NumberFormat
